# Consolidate multi-run title/text-box text into single runs on slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "Title 1" - currently split into many single-word runs that
# concatenate to the same string already, so assigning the identical text
# back would be a no-op. Force a change first, then set the final text so
# the writer rebuilds it as a single consolidated run.
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = " "
$titleShape.TextFrame.TextRange.Text = "A Table, with a caption"

# Shape 3: "TextBox 3" - same situation.
$captionShape = $s.Shapes.Item(3)
$captionShape.TextFrame.TextRange.Text = " "
$captionShape.TextFrame.TextRange.Text = "Demonstration of simple table syntax, with alignment"
